$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_5_4_24"
$ws.Range("B2").Value = 0.5691211202310158
$ws.Range("C2").Value = 0.0902626750975859
$ws.Range("D2").Value = -0.9516738722224041
$ws.Range("E2").Value = -0.5082547403938691
$ws.Range("F2").Value = 0.4768559634685516
$ws.Range("G2").Value = 0.4989897906780243
$ws.Range("H2").Value = 2.413866996765137
$ws.Range("I2").Value = 1.400108218193054

$ws.Range("A3").Value = "model_5_4_23"
$ws.Range("B3").Value = 0.5757217275491305
$ws.Range("C3").Value = 0.1172036977443349
$ws.Range("D3").Value = -0.9117445835005011
$ws.Range("E3").Value = -0.4747926401655296
$ws.Range("F3").Value = 0.4695510566234589
$ws.Range("G3").Value = 0.4842126667499542
$ws.Range("H3").Value = 2.364481687545776
$ws.Range("I3").Value = 1.369045615196228

$ws.Range("A4").Value = "model_5_4_21"
$ws.Range("B4").Value = 0.5825350220256742
$ws.Range("C4").Value = 0.1326685308444797
$ws.Range("D4").Value = -0.8538342242837205
$ws.Range("E4").Value = -0.4336453859804767
$ws.Range("F4").Value = 0.4620107114315033
$ws.Range("G4").Value = 0.4757302701473236
$ws.Range("H4").Value = 2.29285717010498
$ws.Range("I4").Value = 1.330848574638367

$ws.Range("A5").Value = "model_5_4_22"
$ws.Range("B5").Value = 0.5830288571013675
$ws.Range("C5").Value = 0.1449421342644455
$ws.Range("D5").Value = -0.8674389797694404
$ws.Range("E5").Value = -0.4383347481404429
$ws.Range("F5").Value = 0.4614641964435577
$ws.Range("G5").Value = 0.4689981639385223
$ws.Range("H5").Value = 2.309683561325073
$ws.Range("I5").Value = 1.335201740264893

$ws.Range("A6").Value = "model_5_4_20"
$ws.Range("B6").Value = 0.5869596497001563
$ws.Range("C6").Value = 0.1363878043336898
$ws.Range("D6").Value = -0.8148214746598317
$ws.Range("E6").Value = -0.4080214758282488
$ws.Range("F6").Value = 0.4571139812469482
$ws.Range("G6").Value = 0.4736903011798859
$ws.Range("H6").Value = 2.244605302810669
$ws.Range("I6").Value = 1.307061910629272

$ws.Range("A7").Value = "model_5_4_19"
$ws.Range("B7").Value = 0.5947901155764009
$ws.Range("C7").Value = 0.1582087250158815
$ws.Range("D7").Value = -0.7631710446291677
$ws.Range("E7").Value = -0.368812380262975
$ws.Range("F7").Value = 0.4484479129314423
$ws.Range("G7").Value = 0.4617214798927307
$ws.Range("H7").Value = 2.180723190307617
$ws.Range("I7").Value = 1.27066445350647

$ws.Range("A8").Value = "model_5_4_18"
$ws.Range("B8").Value = 0.6030147462379694
$ws.Range("C8").Value = 0.1821708656604507
$ws.Range("D8").Value = -0.7094659176813991
$ws.Range("E8").Value = -0.3276425830763263
$ws.Range("F8").Value = 0.4393456280231476
$ws.Range("G8").Value = 0.4485782384872437
$ws.Range("H8").Value = 2.114299535751343
$ws.Range("I8").Value = 1.232446432113647

$ws.Range("A9").Value = "model_5_4_17"
$ws.Range("B9").Value = 0.6126576245433305
$ws.Range("C9").Value = 0.2125569348453052
$ws.Range("D9").Value = -0.6496164995102263
$ws.Range("E9").Value = -0.2806138056815117
$ws.Range("F9").Value = 0.4286738932132721
$ws.Range("G9").Value = 0.4319115579128265
$ws.Range("H9").Value = 2.040276765823364
$ws.Range("I9").Value = 1.188789844512939

$ws.Range("A10").Value = "model_5_4_16"
$ws.Range("B10").Value = 0.6192047084676988
$ws.Range("C10").Value = 0.224340179852717
$ws.Range("D10").Value = -0.5985727499284479
$ws.Range("E10").Value = -0.2449239740798814
$ws.Range("F10").Value = 0.4214281737804413
$ws.Range("G10").Value = 0.4254484474658966
$ws.Range("H10").Value = 1.977144718170166
$ws.Range("I10").Value = 1.155659198760986

$ws.Range("A11").Value = "model_5_4_15"
$ws.Range("B11").Value = 0.6303485857265281
$ws.Range("C11").Value = 0.2890830269140843
$ws.Range("D11").Value = -0.5444497747617327
$ws.Range("E11").Value = -0.1907363431423805
$ws.Range("F11").Value = 0.4090951979160309
$ws.Range("G11").Value = 0.3899370729923248
$ws.Range("H11").Value = 1.910204648971558
$ws.Range("I11").Value = 1.105356931686401

$ws.Range("A12").Value = "model_5_4_14"
$ws.Range("B12").Value = 0.6411444120900496
$ws.Range("C12").Value = 0.3254591558345579
$ws.Range("D12").Value = -0.4776954694119127
$ws.Range("E12").Value = -0.1375028960051281
$ws.Range("F12").Value = 0.3971473872661591
$ws.Range("G12").Value = 0.3699848353862762
$ws.Range("H12").Value = 1.827641487121582
$ws.Range("I12").Value = 1.055940508842468

$ws.Range("A13").Value = "model_5_4_13"
$ws.Range("B13").Value = 0.6508544813481829
$ws.Range("C13").Value = 0.3538922418692115
$ws.Range("D13").Value = -0.4136120233720559
$ws.Range("E13").Value = -0.08842936462854545
$ws.Range("F13").Value = 0.3864012062549591
$ws.Range("G13").Value = 0.3543893098831177
$ws.Range("H13").Value = 1.748381972312927
$ws.Range("I13").Value = 1.010385632514954

$ws.Range("A14").Value = "model_5_4_12"
$ws.Range("B14").Value = 0.6608365895407475
$ws.Range("C14").Value = 0.3895726738297034
$ws.Range("D14").Value = -0.3505975782731792
$ws.Range("E14").Value = -0.03775926661919771
$ws.Range("F14").Value = 0.3753539621829987
$ws.Range("G14").Value = 0.3348186612129211
$ws.Range("H14").Value = 1.670444488525391
$ws.Range("I14").Value = 0.9633487462997437

$ws.Range("A15").Value = "model_5_4_0"
$ws.Range("B15").Value = 0.6652050574825216
$ws.Range("C15").Value = 0.8094890918983921
$ws.Range("D15").Value = 0.1194445770353357
$ws.Range("E15").Value = 0.3883069765450271
$ws.Range("F15").Value = 0.3705193698406219
$ws.Range("G15").Value = 0.1044950038194656
$ws.Range("H15").Value = 1.08908748626709
$ws.Range("I15").Value = 0.5678327083587646

$ws.Range("A16").Value = "model_5_4_10"
$ws.Range("B16").Value = 0.6661986456046345
$ws.Range("C16").Value = 0.42428848626959
$ws.Range("D16").Value = -0.2952828169303059
$ws.Range("E16").Value = 0.007782698355435658
$ws.Range("F16").Value = 0.3694197535514832
$ws.Range("G16").Value = 0.315777063369751
$ws.Range("H16").Value = 1.602030158042908
$ws.Range("I16").Value = 0.9210723042488098

$ws.Range("A17").Value = "model_5_4_6"
$ws.Range("B17").Value = 0.6686636330058402
$ws.Range("C17").Value = 0.4562410851948439
$ws.Range("D17").Value = -0.202982005835769
$ws.Range("E17").Value = 0.07564846444243745
$ws.Range("F17").Value = 0.3666917383670807
$ws.Range("G17").Value = 0.2982511222362518
$ws.Range("H17").Value = 1.487870812416077
$ws.Range("I17").Value = 0.8580726385116577

$ws.Range("A18").Value = "model_5_4_11"
$ws.Range("B18").Value = 0.6695207233033988
$ws.Range("C18").Value = 0.4246111930404809
$ws.Range("D18").Value = -0.2942770433816164
$ws.Range("E18").Value = 0.008513910361200083
$ws.Range("F18").Value = 0.3657431602478027
$ws.Range("G18").Value = 0.3156000375747681
$ws.Range("H18").Value = 1.600786209106445
$ws.Range("I18").Value = 0.9203935265541077

$ws.Range("A19").Value = "model_5_4_9"
$ws.Range("B19").Value = 0.6773659299799941
$ws.Range("C19").Value = 0.4557261755818756
$ws.Range("D19").Value = -0.2198875873672914
$ws.Range("E19").Value = 0.06488839020745729
$ws.Range("F19").Value = 0.3570608496665955
$ws.Range("G19").Value = 0.29853355884552
$ws.Range("H19").Value = 1.508779883384705
$ws.Range("I19").Value = 0.8680613040924072

$ws.Range("A20").Value = "model_5_4_8"
$ws.Range("B20").Value = 0.6842586257067738
$ws.Range("C20").Value = 0.4724566655946508
$ws.Range("D20").Value = -0.1634156271350462
$ws.Range("E20").Value = 0.1055291941810834
$ws.Range("F20").Value = 0.349432647228241
$ws.Range("G20").Value = 0.2893568873405457
$ws.Range("H20").Value = 1.438934326171875
$ws.Range("I20").Value = 0.830334484577179

$ws.Range("A21").Value = "model_5_4_7"
$ws.Range("B21").Value = 0.687827806031325
$ws.Range("C21").Value = 0.4749818267320272
$ws.Range("D21").Value = -0.1197331389839751
$ws.Range("E21").Value = 0.1337080588279972
$ws.Range("F21").Value = 0.3454826176166534
$ws.Range("G21").Value = 0.2879718244075775
$ws.Range("H21").Value = 1.384907007217407
$ws.Range("I21").Value = 0.8041761517524719

$ws.Range("A22").Value = "model_5_4_5"
$ws.Range("B22").Value = 0.6886388786653624
$ws.Range("C22").Value = 0.529231763275879
$ws.Range("D22").Value = -0.0971935542245661
$ws.Range("E22").Value = 0.1648095237920342
$ws.Range("F22").Value = 0.3445850014686584
$ws.Range("G22").Value = 0.2582157850265503
$ws.Range("H22").Value = 1.357029557228088
$ws.Range("I22").Value = 0.7753047943115234

$ws.Range("A23").Value = "model_5_4_1"
$ws.Range("B23").Value = 0.708778543035911
$ws.Range("C23").Value = 0.8203162773085484
$ws.Range("D23").Value = 0.2296833776102526
$ws.Range("E23").Value = 0.4608131319850249
$ws.Range("F23").Value = 0.3222963511943817
$ws.Range("G23").Value = 0.09855630248785019
$ws.Range("H23").Value = 0.9527420401573181
$ws.Range("I23").Value = 0.5005255341529846

$ws.Range("A24").Value = "model_5_4_2"
$ws.Range("B24").Value = 0.7130113704538446
$ws.Range("C24").Value = 0.8271126296306717
$ws.Range("D24").Value = 0.213306714666282
$ws.Range("E24").Value = 0.4526702973819446
$ws.Range("F24").Value = 0.3176118731498718
$ws.Range("G24").Value = 0.09482851624488831
$ws.Range("H24").Value = 0.9729970693588257
$ws.Range("I24").Value = 0.5080844759941101

$ws.Range("A25").Value = "model_5_4_4"
$ws.Range("B25").Value = 0.7139981687733101
$ws.Range("C25").Value = 0.7296132268684208
$ws.Range("D25").Value = 0.1610661445493614
$ws.Range("E25").Value = 0.389417543632951
$ws.Range("F25").Value = 0.3165197372436523
$ws.Range("G25").Value = 0.14830681681633
$ws.Range("H25").Value = 1.037609100341797
$ws.Range("I25").Value = 0.5668017864227295

$ws.Range("A26").Value = "model_5_4_3"
$ws.Range("B26").Value = 0.7173868464307285
$ws.Range("C26").Value = 0.8107626208024141
$ws.Range("D26").Value = 0.2117193586697858
$ws.Range("E26").Value = 0.4465605234331921
$ws.Range("F26").Value = 0.3127694725990295
$ws.Range("G26").Value = 0.1037964671850204
$ws.Range("H26").Value = 0.9749603271484375
$ws.Range("I26").Value = 0.5137561559677124
